{"js": "const replacements = [\n  [\"2024-07-07 Sunday\", \"2024-07-08 Monday\"],\n  [\"90\u00d735=\", \"76\u00d759=\"],\n  [\"19\u00d713=\", \"52\u00d782=\"],\n  [\"31\u00d769=\", \"16\u00d715=\"],\n  [\"94\u00d723=\", \"59\u00d723=\"],\n  [\"81\u00d725=\", \"21\u00d711=\"],\n  [\"25\u00d773=\", \"46\u00d745=\"],\n  [\"80\u00d769=\", \"78\u00d768=\"],\n  [\"87\u00d729=\", \"20\u00d761=\"],\n  [\"61\u00d755=\", \"31\u00d767=\"],\n  [\"87\u00d787=\", \"11\u00d740=\"],\n  [\"14\u00d773=\", \"69\u00d759=\"],\n  [\"74\u00d745=\", \"19\u00d727=\"],\n  [\"17\u00d790=\", \"84\u00d744=\"],\n  [\"17\u00d730=\", \"73\u00d712=\"],\n  [\"65\u00d711=\", \"99\u00d780=\"],\n  [\"65\u00d765=\", \"48\u00d787=\"],\n  [\"18\u00d749=\", \"48\u00d794=\"],\n  [\"84\u00d788=\", \"93\u00d762=\"],\n  [\"73\u00d736=\", \"12\u00d745=\"],\n  [\"45\u00d726=\", \"78\u00d726=\"],\n  [\"50\u00d792=\", \"13\u00d788=\"],\n  [\"43\u00d759=\", \"59\u00d721=\"],\n  [\"48\u00d770=\", \"42\u00d772=\"],\n  [\"55\u00d733=\", \"40\u00d789=\"],\n  [\"22\u00d713=\", \"31\u00d711=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-07-07 Sunday\", \"2024-07-08 Monday\"),\n    @(\"90\u00d735=\", \"76\u00d759=\"),\n    @(\"19\u00d713=\", \"52\u00d782=\"),\n    @(\"31\u00d769=\", \"16\u00d715=\"),\n    @(\"94\u00d723=\", \"59\u00d723=\"),\n    @(\"81\u00d725=\", \"21\u00d711=\"),\n    @(\"25\u00d773=\", \"46\u00d745=\"),\n    @(\"80\u00d769=\", \"78\u00d768=\"),\n    @(\"87\u00d729=\", \"20\u00d761=\"),\n    @(\"61\u00d755=\", \"31\u00d767=\"),\n    @(\"87\u00d787=\", \"11\u00d740=\"),\n    @(\"14\u00d773=\", \"69\u00d759=\"),\n    @(\"74\u00d745=\", \"19\u00d727=\"),\n    @(\"17\u00d790=\", \"84\u00d744=\"),\n    @(\"17\u00d730=\", \"73\u00d712=\"),\n    @(\"65\u00d711=\", \"99\u00d780=\"),\n    @(\"65\u00d765=\", \"48\u00d787=\"),\n    @(\"18\u00d749=\", \"48\u00d794=\"),\n    @(\"84\u00d788=\", \"93\u00d762=\"),\n    @(\"73\u00d736=\", \"12\u00d745=\"),\n    @(\"45\u00d726=\", \"78\u00d726=\"),\n    @(\"50\u00d792=\", \"13\u00d788=\"),\n    @(\"43\u00d759=\", \"59\u00d721=\"),\n    @(\"48\u00d770=\", \"42\u00d772=\"),\n    @(\"55\u00d733=\", \"40\u00d789=\"),\n    @(\"22\u00d713=\", \"31\u00d711=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        [ref]$oldText,\n        $true,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        [ref]$newText,\n        2\n    ) | Out-Null\n}"}
